$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "9+76="
$t.Cell(1,2).Range.Text = "38+17="
$t.Cell(1,3).Range.Text = "16+59="
$t.Cell(1,4).Range.Text = "73-69="
$t.Cell(1,5).Range.Text = "72-14="
$t.Cell(2,1).Range.Text = "71-3="
$t.Cell(2,2).Range.Text = "27+47="
$t.Cell(2,3).Range.Text = "54-38="
$t.Cell(2,4).Range.Text = "8+4="
$t.Cell(2,5).Range.Text = "5+9="
$t.Cell(3,1).Range.Text = "80-4="
$t.Cell(3,2).Range.Text = "44+17="
$t.Cell(3,3).Range.Text = "64+7="
$t.Cell(3,4).Range.Text = "39+32="
$t.Cell(3,5).Range.Text = "6+58="
$t.Cell(4,1).Range.Text = "25+19="
$t.Cell(4,2).Range.Text = "6+35="
$t.Cell(4,3).Range.Text = "29+17="
$t.Cell(4,4).Range.Text = "17+34="
$t.Cell(4,5).Range.Text = "72-25="
$t.Cell(5,1).Range.Text = "26+28="
$t.Cell(5,2).Range.Text = "11-7="
$t.Cell(5,3).Range.Text = "8+49="
$t.Cell(5,4).Range.Text = "55+16="
$t.Cell(5,5).Range.Text = "46+26="
$t.Cell(6,1).Range.Text = "41-4="
$t.Cell(6,2).Range.Text = "17+38="
$t.Cell(6,3).Range.Text = "55+39="
$t.Cell(6,4).Range.Text = "38+15="
$t.Cell(6,5).Range.Text = "73-67="
$t.Cell(7,1).Range.Text = "26+68="
$t.Cell(7,2).Range.Text = "63-49="
$t.Cell(7,3).Range.Text = "9+43="
$t.Cell(7,4).Range.Text = "92-44="
$t.Cell(7,5).Range.Text = "86+6="
$t.Cell(8,1).Range.Text = "46+7="
$t.Cell(8,2).Range.Text = "15+26="
$t.Cell(8,3).Range.Text = "74-59="
$t.Cell(8,4).Range.Text = "5+66="
$t.Cell(8,5).Range.Text = "63-9="
$t.Cell(9,1).Range.Text = "9+9="
$t.Cell(9,2).Range.Text = "53+19="
$t.Cell(9,3).Range.Text = "69+19="
$t.Cell(9,4).Range.Text = "76-59="
$t.Cell(9,5).Range.Text = "77+18="
$t.Cell(10,1).Range.Text = "33-25="
$t.Cell(10,2).Range.Text = "83-75="
$t.Cell(10,3).Range.Text = "95-77="
$t.Cell(10,4).Range.Text = "91-52="
$t.Cell(10,5).Range.Text = "26+25="
$t.Cell(11,1).Range.Text = "94-16="
$t.Cell(11,2).Range.Text = "92-86="
$t.Cell(11,3).Range.Text = "27+36="
$t.Cell(11,4).Range.Text = "9+54="
$t.Cell(11,5).Range.Text = "54+9="
$t.Cell(12,1).Range.Text = "60-46="
$t.Cell(12,2).Range.Text = "74-25="
$t.Cell(12,3).Range.Text = "94-57="
$t.Cell(12,4).Range.Text = "5+18="
$t.Cell(12,5).Range.Text = "76-58="
$t.Cell(13,1).Range.Text = "36+28="
$t.Cell(13,2).Range.Text = "71-17="
$t.Cell(13,3).Range.Text = "26+5="
$t.Cell(13,4).Range.Text = "59+32="
$t.Cell(13,5).Range.Text = "90-85="
$t.Cell(14,1).Range.Text = "80-16="
$t.Cell(14,2).Range.Text = "66-38="
$t.Cell(14,3).Range.Text = "36+25="
$t.Cell(14,4).Range.Text = "68+13="
$t.Cell(14,5).Range.Text = "54+28="
$t.Cell(15,1).Range.Text = "51-49="
$t.Cell(15,2).Range.Text = "26+28="
$t.Cell(15,3).Range.Text = "15-8="
$t.Cell(15,4).Range.Text = "4+47="
$t.Cell(15,5).Range.Text = "50-39="
$t.Cell(16,1).Range.Text = "29+56="
$t.Cell(16,2).Range.Text = "43-5="
$t.Cell(16,3).Range.Text = "24+39="
$t.Cell(16,4).Range.Text = "19+65="
$t.Cell(16,5).Range.Text = "76-27="
$t.Cell(17,1).Range.Text = "70-29="
$t.Cell(17,2).Range.Text = "44+38="
$t.Cell(17,3).Range.Text = "72-39="
$t.Cell(17,4).Range.Text = "13+38="
$t.Cell(17,5).Range.Text = "61-53="
$t.Cell(18,1).Range.Text = "81-4="
$t.Cell(18,2).Range.Text = "17+65="
$t.Cell(18,3).Range.Text = "75-29="
$t.Cell(18,4).Range.Text = "67+27="
$t.Cell(18,5).Range.Text = "91-45="
$t.Cell(19,1).Range.Text = "83-15="
$t.Cell(19,2).Range.Text = "15+19="
$t.Cell(19,3).Range.Text = "56+16="
$t.Cell(19,4).Range.Text = "92-87="
$t.Cell(19,5).Range.Text = "64-35="
$t.Cell(20,1).Range.Text = "90-4="
$t.Cell(20,2).Range.Text = "95-38="
$t.Cell(20,3).Range.Text = "62-59="
$t.Cell(20,4).Range.Text = "9+7="
$t.Cell(20,5).Range.Text = "84-57="
